# LNI-317: Added previously removed signature blocks back to test Statutory Instruments.
#
# Appends (after the existing final paragraph) six new paragraphs:
#   1. empty N3 paragraph with numbering removed
#   2. empty SigBlock paragraph
#   3. SigBlock paragraph: tab + "Senior Official" (character style SigSignee)
#   4. SigBlock paragraph: tab + "A senior officer of the" (character style Sigtitle)
#   5. SigBlock paragraph: tab + "Department of Agriculture, Environment and Rural Affairs"
#      (character style Sigtitle)
#   6. empty N3 paragraph with numbering removed

$d = $word.ActiveDocument

# Collapse to the very end of the document body and insert the new
# paragraphs there (i.e. immediately before the final section break).
$endRange = $d.Content
$endRange.Collapse(0)

$sel = $word.Selection
$sel.SetRange($endRange.Start, $endRange.End)

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>Senior Official</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>A senior officer of the</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>Department of Agriculture, Environment and Rural Affairs</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr></w:p>
'@

[void]$sel.InsertXML($newParagraphsXml)

# Apply the character styles to the signatory lines via Find & Replace
# (replacing the text with itself while stamping the Replacement.Style)
# so the runs pick up the correct w:rStyle without disturbing anything
# else in the freshly-inserted paragraphs.
function Set-RunStyle($text, $styleName) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $text
    $find.Replacement.Text = $text
    $find.Replacement.Style = $styleName
    [void]$find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

Set-RunStyle "Senior Official" "SigSignee"
Set-RunStyle "A senior officer of the" "Sigtitle"
Set-RunStyle "Department of Agriculture, Environment and Rural Affairs" "Sigtitle"

Write-Host "Signature block appended."
